$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column B ("Room Helper") before the existing "Seats" column,
#    shifting the old column B ("Seats") to column C.
$ws.Columns("B").Insert()

# 2. Header row.
$ws.Range("B1").Value2 = "Room Helper"

# 3. Append the four new rooms reported by the sponsor to the bottom of the
#    table (rows 39-42) before filling in the helper formula column so the
#    shared formula can cover the full B2:B42 range in one shot.
$ws.Range("A39").Value2 = 117
$ws.Range("C39").Value2 = 40

$ws.Range("A40").Value2 = 127
$ws.Range("C40").Value2 = 20

$ws.Range("A41").Value2 = 130
$ws.Range("C41").Value2 = 38

# Row 42's room number was entered via a TEXT() formula (matching the
# original author's file) rather than as a literal number.
$ws.Range("A42").Formula = "=TEXT(248,0)"
$ws.Range("C42").Value2 = 40

# 4. Fill the new "Room Helper" column with a TEXT() formula that mirrors
#    column A as text, for every data row (2 through 42). B2 is entered on
#    its own (standalone formula) and B3:B42 is entered as one range so
#    Excel collapses that block into a single shared formula, matching the
#    original author's edit.
$ws.Range("B2").Formula = "=TEXT(A2,0)"
$ws.Range("B3:B42").Formula = "=TEXT(A3,0)"

# 5. Restore the view state recorded in the edited workbook: scrolled so row
#    8 is the top-left visible row, with B2:B42 selected (active cell B2).
$win = $excel.ActiveWindow
$win.ScrollRow = 8
$win.ScrollColumn = 1
$ws.Range("B2:B42").Select()
